# Hortaliza, Macroferia Regional de Talca - Betarraga
# A new weekly price record is inserted at the top of the data block
# (row 564), pushing all subsequent records down by one row and
# extending the used range from R689 to R690.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 564; this shifts rows
# 564:689 down to 565:690 (and their contents) automatically.
$ws.Rows.Item(564).Insert()

# Populate the new row 564 with the new weekly record.
$ws.Range("A564").Value = 5
$ws.Range("B564").Value = "Macroferia Regional de Talca"
$ws.Range("C564").Value = "Maule"
$ws.Range("D564").Value = 45244
$ws.Range("E564").Value = 7
$ws.Range("F564").Value = 100114014
$ws.Range("G564").Value = "Betarraga"
$ws.Range("H564").Value = "Sin especificar"
$ws.Range("I564").Value = "Primera"
$ws.Range("J564").Value = 4000
$ws.Range("K564").Value = 800
$ws.Range("L564").Value = 800
$ws.Range("M564").Value = 800
$ws.Range("N564").Value = "`$/paquete 5 unidades"
$ws.Range("O564").Value = "Región del Maule"
$ws.Range("P564").Value = 160
$ws.Range("Q564").Value = 5
$ws.Range("R564").Value = "Hortaliza"
